$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Category"
$ws.Range("B1").Value = "Words included"

$ws.Range("B2").Select()
